# Generate data from FS Sugeno
# Insert 6 new "V_*" input columns before the existing CLPVariation column,
# then drop the old Critical / FinalOut columns, leaving:
# A..F (unchanged) | G..L (V_*) | M CLPVariation | N CLPVariation_pred | O erro_CLP

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert six blank columns at G..L (shifts old G:K -> M:Q)
$ws.Range("G1:L1").EntireColumn.Insert()

# 2) Remove the old "Critical" / "FinalOut" columns, now at O:P
#    (old layout after insert: M=CLPVariation, N=CLPVariation_pred, O=Critical, P=FinalOut, Q=erro_CLP)
$ws.Range("O1:P1").EntireColumn.Delete()

# 3) Header row for the new V_* columns
$ws.Range("G1").Value = "V_MemoryUsage"
$ws.Range("H1").Value = "V_ProcessorLoad"
$ws.Range("I1").Value = "V_InpNetThroughput"
$ws.Range("J1").Value = "V_OutNetThroughput"
$ws.Range("K1").Value = "V_OutBandwidth"
$ws.Range("L1").Value = "V_Latency"

# 4) Data rows for the new V_* columns
$ws.Range("G2").Value = 0.45
$ws.Range("H2").Value = 0.52
$ws.Range("I2").Value = 0.5
$ws.Range("J2").Value = 0.48
$ws.Range("K2").Value = 0.45
$ws.Range("L2").Value = 0.48

$ws.Range("G3").Value = 0.6
$ws.Range("H3").Value = 0.52
$ws.Range("I3").Value = 0.5
$ws.Range("J3").Value = 0.52
$ws.Range("K3").Value = 0.48
$ws.Range("L3").Value = 0.52

$ws.Range("G4").Value = 0.45
$ws.Range("H4").Value = 0.52
$ws.Range("I4").Value = 0.5
$ws.Range("J4").Value = 0.48
$ws.Range("K4").Value = 0.46
$ws.Range("L4").Value = 0.48

$ws.Range("G5").Value = 0.5
$ws.Range("H5").Value = 0.5
$ws.Range("I5").Value = 0.5
$ws.Range("J5").Value = 0.5
$ws.Range("K5").Value = 0.5
$ws.Range("L5").Value = 0.5

$ws.Range("G6").Value = 0.5
$ws.Range("H6").Value = 0.5
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.5
$ws.Range("K6").Value = 0.5
$ws.Range("L6").Value = 0.5

$ws.Range("G7").Value = 0.4
$ws.Range("H7").Value = 0.45
$ws.Range("I7").Value = 0.51
$ws.Range("J7").Value = 0.45
$ws.Range("K7").Value = 0.3
$ws.Range("L7").Value = 0.55

$ws.Range("G8").Value = 0.59
$ws.Range("H8").Value = 0.45
$ws.Range("I8").Value = 0.55
$ws.Range("J8").Value = 0.52
$ws.Range("K8").Value = 0.49
$ws.Range("L8").Value = 0.52

$ws.Range("G9").Value = 0.45
$ws.Range("H9").Value = 0.53
$ws.Range("I9").Value = 0.56
$ws.Range("J9").Value = 0.48
$ws.Range("K9").Value = 0.46
$ws.Range("L9").Value = 0.48

$ws.Range("G10").Value = 0.56
$ws.Range("H10").Value = 0.67
$ws.Range("I10").Value = 0.58
$ws.Range("J10").Value = 0.55
$ws.Range("K10").Value = 0.52
$ws.Range("L10").Value = 0.6

$ws.Range("G11").Value = 0.56
$ws.Range("H11").Value = 0.67
$ws.Range("I11").Value = 0.58
$ws.Range("J11").Value = 0.55
$ws.Range("K11").Value = 0.52
$ws.Range("L11").Value = 0.6
